# Applies the "12.Woche" report entry to Wochenberichte.docx
# (Kommentierung des Codes und Implementierung der Email-Pushmitteilungen in GPIO_Test)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 0. The trailing "_GoBack" bookmark currently sits at the end of the
#    last written paragraph ("... erfolgreich erstellt."). It needs to
#    move to the end of the new last paragraph, so drop it here and
#    re-create it once the new content exists.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# The document already ends with a single empty paragraph right before
# the section break - keep it as the blank line that separates the
# previous entry from the new heading, and append everything else
# after it, in plain/unformatted text first so later character
# formatting on the heading cannot "leak" into the other paragraphs.
$anchor = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$anchor.Collapse(0)

# 1. Heading paragraph text: "12.Woche 23.05.2018 - 30.05.2018:"
$headingPara = $d.Paragraphs.Add($anchor)
$headingRange = $headingPara.Range
$headingRange.Text = "12.Woche 23.05.2018 " + [char]0x2013 + " 30.05.2018:"

# 2. Blank separator line
$anchor2 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$anchor2.Collapse(0)
$blankPara = $d.Paragraphs.Add($anchor2)

# 3. "Martin Eller hat die E-Mail Pushmitteilungen in die GPIO_Test ..."
$anchor3 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$anchor3.Collapse(0)
$martinPara = $d.Paragraphs.Add($anchor3)
$martinRange = $martinPara.Range
$martinRange.Text = "Martin Eller hat die E-Mail Pushmitteilungen in die GPIO_Test" + `
    " Klasse implementiert. Es werden nun erfolgreich beim Schließen bzw. " + `
    "Öffnen des Postkastens Emails verschickt, welche nebenbei dann auch " + `
    "noch das Datum und die Uhrzeit der letzten Schließung erhalten. "

# 4. "Ebenfalls wurde gewährleistet, ..."
$anchor4 = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$anchor4.Collapse(0)
$finalPara = $d.Paragraphs.Add($anchor4)
$finalRange = $finalPara.Range
$finalRange.Text = "Ebenfalls wurde gewährleistet, dass das Programm gleich nach " + `
    "dem Starten des Raspberry von einen Skript ausgeführt wird."

# ---------------------------------------------------------------------
# 5. Re-create the "_GoBack" bookmark at the end of the new last
#    paragraph (collapsed, zero-length, just like the original).
# ---------------------------------------------------------------------
$bookmarkRange = $finalPara.Range
$bookmarkRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------
# 6. Only now apply the heading-specific (Heading 1 - like) character
#    and paragraph formatting, retroactively, to the heading paragraph.
#    Doing this last keeps it from bleeding into the paragraphs created
#    afterwards.
# ---------------------------------------------------------------------
$headingFormat = $headingPara.Range.ParagraphFormat
$headingFormat.KeepWithNext = 1
$headingFormat.KeepTogether = 1
$headingFormat.OutlineLevel = 1
$headingFormat.SpaceBefore = 12
$headingFormat.SpaceAfter = 0

$headingFont = $headingPara.Range.Font
$headingFont.NameAscii = "Calibri Light"
$headingFont.NameFarEast = "Calibri Light"
$headingFont.NameBi = "Calibri Light"
$headingFont.Size = 16
$headingFont.SizeBi = 16
$headingFont.Underline = 1
$headingFont.TextColor.ObjectThemeColor = 4

Write-Host "Week 12 report entry inserted."
